$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D2").Value = "2016-02-22 17:39:23"
$wsZh.Range("G2").Value = "2016-02-22 17:40:10"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D2").Value = "2016-02-22 17:39:35"
$wsDe.Range("G2").Value = "2016-02-22 17:40:31"
